# Third Commit (Ruli) - katalon bio farma v.01
#
# Appends "-AutomatedTest" to each currency code in column C of the
# "Currency" sheet (rows 2-6), and moves the active selection to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Currency")

$ws.Range("C2").Value = "IDR-AutomatedTest"
$ws.Range("C3").Value = "USD-AutomatedTest"
$ws.Range("C4").Value = "EUR-AutomatedTest"
$ws.Range("C5").Value = "EGP-AutomatedTest"
$ws.Range("C6").Value = "KPW-AutomatedTest"

[void]$ws.Range("C10").Select()
